$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format first so numeric-looking price strings
# (e.g. "0.7094", "242.03") are stored verbatim as text instead of being
# auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 42 and 43 swap coin/link content (TrustWalletToken <-> Quant) with
# each also getting an updated price/volume figure.
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '110.43'
$ws.Range("E42").Value = '  +5.41%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '0.8916'
$ws.Range("E43").Value = '  +0.33%  '

# Per-row price / 1h-volume refreshes.
$ws.Range("D2").Value = '29.312.09'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").Value = '1.873.62'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '0.7094'
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").Value = '242.03'
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.07868'
$ws.Range("E8").Value = '  +2.24%  '
$ws.Range("D9").Value = '0.3112'
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("D10").Value = '25.21'
$ws.Range("E10").Value = '  +0.25%  '
$ws.Range("D11").Value = '0.08397'
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").Value = '1.872.33'
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").Value = '5.237'
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("D14").Value = '0.7184'
$ws.Range("E14").Value = '  +1.01%  '
$ws.Range("D15").Value = '91.28'
$ws.Range("E15").Value = '  -0.01%  '
$ws.Range("D16").Value = '0.000008374'
$ws.Range("E16").Value = '  +1.39%  '
$ws.Range("D17").Value = '6.139'
$ws.Range("E17").Value = '  +3.39%  '
$ws.Range("D18").Value = '29.311.74'
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("D19").Value = '240.93'
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("D21").Value = '2.125.50'
$ws.Range("E21").Value = '  -0.72%  '
$ws.Range("D22").Value = '0.9997'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = '7.766'
$ws.Range("E23").Value = '  -1.11%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").Value = '0.1595'
$ws.Range("E25").Value = '  -2.12%  '
$ws.Range("D26").Value = '162.83'
$ws.Range("E26").Value = '  -0.26%  '
$ws.Range("D27").Value = '9.047'
$ws.Range("E27").Value = '  +0.39%  '
$ws.Range("D28").Value = '18.54'
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("D29").Value = '1.507'
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("D30").Value = '4.410'
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").Value = '4.339'
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("E32").Value = '  -4.84%  '
$ws.Range("D33").Value = '0.05356'
$ws.Range("E33").Value = '  +2.09%  '
$ws.Range("E34").Value = '  +0.92%  '
$ws.Range("D35").Value = '1.176'
$ws.Range("E35").Value = '  +0.34%  '
$ws.Range("D36").Value = '0.7473'
$ws.Range("E36").Value = '  -0.91%  '
$ws.Range("D37").Value = '2.684'
$ws.Range("E37").Value = '  +0.09%  '
$ws.Range("E38").Value = '  +1.12%  '
$ws.Range("D39").Value = '1.246.46'
$ws.Range("E39").Value = '  +8.13%  '
$ws.Range("E40").Value = '  +0.58%  '
$ws.Range("E41").Value = '  +2.30%  '
$ws.Range("D44").Value = '72.45'
$ws.Range("E44").Value = '  -0.89%  '
$ws.Range("E45").Value = '  +13.36%  '
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").Value = '2.017.53'
$ws.Range("E47").Value = '  -0.59%  '
$ws.Range("D48").Value = '1.798'
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("D49").Value = '0.5186'
$ws.Range("E49").Value = '  -0.17%  '
$ws.Range("D50").Value = '9.451'
$ws.Range("E50").Value = '  +0.72%  '
$ws.Range("D51").Value = '0.4358'
$ws.Range("E51").Value = '  +1.47%  '

# Strip the temporary text-number-format so the touched cells end up
# unstyled again, matching the rest of the sheet.
$ws.Range("D2:D51").ClearFormats()
